$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(39.78, 0.16, 10.01247477531433, 9.3125),
    @(39.78, 0.16, 9.666251182556152, 9.03125),
    @(39.78, 0.16, 9.286051273345947, 8.796875),
    @(39.78, 0.16, 9.419345617294312, 8.84375)
)

$startRow = 108
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
    $ws.Cells.Item($row, 4).Value = $data[$i][3]
}
